$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.379.26"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "2.440.60"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.48"
$ws.Range("E5").Value = "  -1.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.84"
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.527"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.109"
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.14"
$ws.Range("E11").Value = "  -2.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.345"
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.37"
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.892.75"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000172"
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("D16").Value = "62.454.30"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").Value = "2.443.04"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.65"
$ws.Range("E18").Value = "  -6.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.65"
$ws.Range("E19").Value = "  -3.78%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "319.01"
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.10"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.18"
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.72"
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.59"
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "638.02"
$ws.Range("E26").Value = "  -4.99%  "
$ws.Range("D27").Value = "2.574.56"
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0940"
$ws.Range("E29").Value = "  -5.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.39"
$ws.Range("E30").Value = "  -4.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.73"
$ws.Range("E31").Value = "  -4.28%  "
$ws.Range("E32").Value = "  -4.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.130"
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.47"
$ws.Range("E35").Value = "  -4.86%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "151.86"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.58"
$ws.Range("E37").Value = "  -4.62%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.361"
$ws.Range("E38").Value = "  -2.84%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.40"
$ws.Range("E39").Value = "  -2.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.22"
$ws.Range("E40").Value = "  -4.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.62"
$ws.Range("E41").Value = "  -3.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.68"
$ws.Range("E42").Value = "  -4.08%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").Value = "0.0₆0306"
$ws.Range("E44").Value = "  -1.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "151.76"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.39"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.49"
$ws.Range("E47").Value = "  -3.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.599"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.71"
$ws.Range("E49").Value = "  -5.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0498"
$ws.Range("E50").Value = "  -3.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0897"
$ws.Range("E51").Value = "  -2.55%  "
